# Mini Project Class One
# Color the entire row for SNO 11 (Custom Services / Development of Custom
# Services / 2 hrs / 1 day) green (RGB 00B050) to match wdColorGreen-style
# highlight used elsewhere in the syllabus (e.g. row 9's "2 days").

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$row = $t.Rows.Item(15)

# wdColor value for RGB hex 00B050 (R=0x00,G=0xB0,B=0x50) stored as BGR ->
# 0x0050B000 in Word's color space = 0x50B000 = 5287936
$green = 5287936

for ($i = 1; $i -le $row.Cells.Count; $i++) {
    $cell = $row.Cells.Item($i)
    $cellRange = $cell.Range
    $cellRange.Font.Color = $green
}
